# 自动更新Excel文件 - 2025-11-30 23:12:05
#
# For every data row, column F holds a "start date" stored as a plain
# YYYYMMDD integer and column D holds the total countdown length (days).
# Column E ("剩余" / remaining) is recomputed against "today"
# (2025-12-01): remaining = (start date + D days) - today.
# If that would be zero or negative (the countdown already finished),
# the cycle restarts: E is reset back to the full D, and F (the start
# date) is reset to today.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Howard Hinnant's civil_from_days / days_from_civil algorithm, done with
# plain integer math so it is independent of any COM/.NET DateTime support.
function DaysFromCivil([int]$y, [int]$m, [int]$d) {
    if ($m -le 2) { $y = $y - 1 }
    $era = [math]::Floor($y / 400)
    $yoe = $y - $era * 400
    $mp = ($m + 9) % 12
    $doy = [math]::Floor((153 * $mp + 2) / 5) + $d - 1
    $doe = $yoe * 365 + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100) + $doy
    return $era * 146097 + $doe - 719468
}

$todayYear = 2025
$todayMonth = 12
$todayDay = 1
$todaySerial = DaysFromCivil $todayYear $todayMonth $todayDay
$todayNum = 20251201

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $totalDays = $dCell.Value()
    $fRaw = $fCell.Value()

    if ($totalDays -eq $null -or $fRaw -eq $null) {
        continue
    }

    $fStr = [string]([int64]$fRaw)

    # Skip malformed / unparsable start dates (not exactly 8 digits,
    # e.g. a stray "202510929") - leave the row untouched.
    if ($fStr.Length -ne 8) {
        continue
    }

    $startYear = [int]$fStr.Substring(0, 4)
    $startMonth = [int]$fStr.Substring(4, 2)
    $startDay = [int]$fStr.Substring(6, 2)

    $startSerial = DaysFromCivil $startYear $startMonth $startDay
    $endSerial = $startSerial + [int]$totalDays

    $remaining = $endSerial - $todaySerial

    if ($remaining -le 0) {
        # Countdown finished - restart the cycle from today.
        $eCell.Value = $totalDays
        $fCell.Value = $todayNum
    } else {
        $eCell.Value = $remaining
    }
}
